$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) column F
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 4525
$wsExhibition.Range("F3").Value = 860
$wsExhibition.Range("F7").Value = 160
$wsExhibition.Range("F8").Value = 640
$wsExhibition.Range("F10").Value = 204
$wsExhibition.Range("F11").Value = 1420
$wsExhibition.Range("F12").Value = 32
$wsExhibition.Range("F13").Value = 3015
$wsExhibition.Range("F14").Value = 450
$wsExhibition.Range("F15").Value = 684

# Sheet "全部类型" (All types) - same events, rows offset by +1 from row 10 onward
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4525
$wsAll.Range("F3").Value = 860
$wsAll.Range("F7").Value = 160
$wsAll.Range("F8").Value = 640
$wsAll.Range("F11").Value = 204
$wsAll.Range("F12").Value = 1420
$wsAll.Range("F13").Value = 32
$wsAll.Range("F14").Value = 3015
$wsAll.Range("F15").Value = 450
$wsAll.Range("F16").Value = 684
